# Weekly update: insert a new daily price record for Limón (Agrícola del
# Norte S.A. de Arica) as row 398, pushing the existing rows 398-426 down
# to 399-427.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 398 (existing rows 398:426 shift to 399:427)
$ws.Rows("398").Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(398, 1).Value  = 1
$ws.Cells.Item(398, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(398, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(398, 4).Value  = 45075
$ws.Cells.Item(398, 5).Value  = 15
$ws.Cells.Item(398, 6).Value  = "Fruta"
$ws.Cells.Item(398, 7).Value  = 100102
$ws.Cells.Item(398, 8).Value  = "Cítricos"
$ws.Cells.Item(398, 9).Value  = 100102003
$ws.Cells.Item(398, 10).Value = "Limón"
$ws.Cells.Item(398, 11).Value = "Sutil De Gase"
$ws.Cells.Item(398, 12).Value = "Primera"
$ws.Cells.Item(398, 13).Value = 250
$ws.Cells.Item(398, 14).Value = 26000
$ws.Cells.Item(398, 15).Value = 27000
$ws.Cells.Item(398, 16).Value = 26500
$ws.Cells.Item(398, 17).Value = "$/caja 24 kilos"
$ws.Cells.Item(398, 18).Value = "Perú"
$ws.Cells.Item(398, 19).Value = 1104
$ws.Cells.Item(398, 20).Value = 24
